# docx writer: avoid extra copy of abstractNum and num elements in
# numbering.xml. Pandoc's reference.docx carried numbering definitions
# that got duplicated into the generated docx (clashing autogenerated
# abstractNum/num elements under the same ids), which confused Word
# Online even though Desktop Word tolerated it. Strip the duplicate
# <w:abstractNum w:abstractNumId="990"> block and the duplicate
# <w:num w:numId="1000"> element, keeping a single copy of each.

$d = $word.ActiveDocument
$xml = $d.WordOpenXML

$openTag = '<w:abstractNum w:abstractNumId="990">'
$closeTag = '</w:abstractNum>'

$first = $xml.IndexOf($openTag)
if ($first -lt 0) {
    throw "could not find abstractNum 990 in WordOpenXML"
}
$second = $xml.IndexOf($openTag, $first + $openTag.Length)
if ($second -lt 0) {
    throw "only one abstractNum 990 found; nothing to dedupe"
}
$secondClose = $xml.IndexOf($closeTag, $second) + $closeTag.Length

# Drop the first (duplicate/clashing) <w:abstractNum w:abstractNumId="990">
# ... </w:abstractNum> block, keeping the second one intact.
$xml = $xml.Substring(0, $first) + $xml.Substring($second, $secondClose - $second) + $xml.Substring($secondClose)

$numTag = '<w:num w:numId="1000"><w:abstractNumId w:val="990" /></w:num>'
$firstNum = $xml.IndexOf($numTag)
if ($firstNum -lt 0) {
    throw "could not find num 1000 in WordOpenXML"
}
$secondNum = $xml.IndexOf($numTag, $firstNum + $numTag.Length)
if ($secondNum -lt 0) {
    throw "only one num 1000 found; nothing to dedupe"
}

# Drop the duplicate <w:num w:numId="1000"> element, keeping one copy.
$xml = $xml.Remove($firstNum, $numTag.Length)

$d.WordOpenXML = $xml
